# add LTW and GOV NRW 2022. also fix Zentrum partyfacts ID and MR Code switchup
#
# Appends 5 new party rows (13-17) to the "Tabelle1" sheet:
#   Sport / Deutsche Sportpartei
#   Liebe / Europäische Partei Liebe
#   neo   / neo. Wohlstand für alle
#   PdF   / Partei des Fortschritts
#   LfK   / Partei für Kinder, Jugendliche und Familien - Lobbyisten für Kinder
# each with its info-url in column H.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13: Deutsche Sportpartei -----------------------------------------
$ws.Range("A13").Value = "Sport"
$ws.Range("B13").Value = "Deutsche Sportpartei"
$ws.Range("H13").Value = "https://vereins.fandom.com/wiki/Deutsche_Sportpartei"

# --- Row 14: Europäische Partei Liebe --------------------------------------
$ws.Range("A14").Value = "Liebe"
$ws.Range("B14").Value = "Europäische Partei Liebe"
$ws.Range("H14").Value = "https://de.wikipedia.org/wiki/Europ%C3%A4ische_Partei_Liebe"
$ws.Range("H14").Interior.ColorIndex = -4142

# --- Row 15: neo. Wohlstand für alle ---------------------------------------
$ws.Range("A15").Value = "neo"
$ws.Range("B15").Value = "neo. Wohlstand für alle"
$ws.Range("H15").Value = "https://www.bpb.de/themen/parteien/wer-steht-zur-wahl/nordrhein-westfalen-2022/507345/neo-wohlstand-fuer-alle/"
$ws.Range("H15").Interior.ColorIndex = -4142

# --- Row 16: Partei des Fortschritts ----------------------------------------
$ws.Range("A16").Value = "PdF"
$ws.Range("B16").Value = "Partei des Fortschritts"
$ws.Range("H16").Value = "https://de.wikipedia.org/wiki/Partei_des_Fortschritts"
$ws.Range("H16").Interior.ColorIndex = -4142

# --- Row 17: Lobbyisten für Kinder -------------------------------------------
$ws.Range("A17").Value = "LfK"
$ws.Range("B17").Value = "Partei für Kinder, Jugendliche und Familien - Lobbyisten für Kinder"
$ws.Range("H17").Value = "https://de.wikipedia.org/wiki/Lobbyisten_f%C3%BCr_Kinder"
$ws.Range("H17").Interior.ColorIndex = -4142

# Leave the selection where the author's last click landed.
$ws.Range("G31").Select() | Out-Null
